$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.677027463912964
$ws.Range("B1").Value = 2.253070116043091
$ws.Range("C1").Value = 1.442201495170593
$ws.Range("D1").Value = 1.514088869094849
$ws.Range("E1").Value = 1.517507195472717
